$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$d = [DateTime]::ParseExact("2021-04-14", "yyyy-MM-dd", $null)

# ---------------------------------------------------------------------
# New defect rows (4-7). Cells are written in the exact order needed so
# newly-introduced shared strings land at the same indices as the target
# workbook (Design, DCM_APP/SetManualFlag, N/A, UTS_APP_TC_2, Fail due
# lack of input validation., UTS_APP_TC_4, UTS_SM_TC_4,
# SerialManager/DBG_Error, UTS_SM_TC_3, SerialManager/DBG_Print,
# DCM_APP/SetTestFlag).
# ---------------------------------------------------------------------
$ws1.Range("F4").Value = "Design"
$ws1.Range("B4").Value = "DCM_APP/SetManualFlag"
$ws1.Range("G4").Value = "N/A"
$ws1.Range("D4").Value = "UTS_APP_TC_2"
$ws1.Range("I4").Value = "Fail due lack of input validation."

$ws1.Range("D5").Value = "UTS_APP_TC_4"
$ws1.Range("D7").Value = "UTS_SM_TC_4"
$ws1.Range("B7").Value = "SerialManager/DBG_Error"
$ws1.Range("D6").Value = "UTS_SM_TC_3"
$ws1.Range("B6").Value = "SerialManager/DBG_Print"
$ws1.Range("B5").Value = "DCM_APP/SetTestFlag"

# --- Row 4 remaining cells ---------------------------------------------
$ws1.Range("C4").Value = $d
$ws1.Range("E4").Value = "Checking"
$ws1.Range("E4").WrapText = $true
$ws1.Range("E4").Font.Color = 0
$ws1.Range("H4").Value = 0

# --- Row 5 remaining cells ---------------------------------------------
$ws1.Range("C5").Value = $d
$ws1.Range("E5").Value = "Checking"
$ws1.Range("E5").WrapText = $true
$ws1.Range("E5").Font.Color = 0
$ws1.Range("F5").Value = "Design"
$ws1.Range("G5").Value = "N/A"
$ws1.Range("H5").Value = 0
$ws1.Range("I5").Value = "Fail due lack of input validation."

# --- Row 6 remaining cells ---------------------------------------------
$ws1.Range("C6").Value = $d
$ws1.Range("E6").Value = "Checking"
$ws1.Range("E6").WrapText = $true
$ws1.Range("E6").Font.Color = 0
$ws1.Range("F6").Value = "Design"
$ws1.Range("G6").Value = "N/A"
$ws1.Range("H6").Value = 0
$ws1.Range("I6").Value = "Fail due lack of input validation."

# --- Row 7 remaining cells ---------------------------------------------
$ws1.Range("C7").Value = $d
$ws1.Range("E7").Value = "Checking"
$ws1.Range("E7").WrapText = $true
$ws1.Range("E7").Font.Color = 0
$ws1.Range("F7").Value = "Design"
$ws1.Range("G7").Value = "N/A"
$ws1.Range("H7").Value = 0
$ws1.Range("I7").Value = "Fail due lack of input validation."

# ---------------------------------------------------------------------
# Column width tweaks: column B widened (manual resize) and column D
# auto-fit to the new, longer test-case identifiers it now holds.
# ---------------------------------------------------------------------
$ws1.Columns("B").ColumnWidth = 25.42578125
$ws1.Columns("D").AutoFit()

# ---------------------------------------------------------------------
# Selection bookkeeping: the author ended up on sheet 2 (selecting B10)
# before returning to sheet 1 (selecting D11), which stays the active tab.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B10").Select()

$ws1.Activate()
$ws1.Range("D11").Select()
